# Apply "slight tweaks in the car times" edit to Rotterdam_car_times.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value tweaks (car times matrix) ---
$ws.Range("M4").Value = 22
$ws.Range("W4").Value = 21
$ws.Range("H9").Value = 22
$ws.Range("Y10").Value = 10
$ws.Range("Y12").Value = 10
$ws.Range("H19").Value = 21
$ws.Range("N21").Value = 10
$ws.Range("P21").Value = 10

# --- Column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 25.90625
$ws.Columns.Item(12).ColumnWidth = 15.36328125
$ws.Columns.Item(13).ColumnWidth = 14.36328125
$ws.Columns.Item(14).ColumnWidth = 19.7265625
$ws.Columns.Item(15).ColumnWidth = 18.1796875
$ws.Columns.Item(16).ColumnWidth = 14.6328125
$ws.Columns.Item(23).ColumnWidth = 13.26953125


# --- View/selection tweaks ---
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("P23").Select()
